# Apply "Further minor improvements and labelling" edits
$wb = $excel.ActiveWorkbook

# --- Rename second sheet: "Timesatser_budget" -> "Timesatser" ---
$wsRates = $wb.Worksheets.Item("Timesatser_budget")
$wsRates.Name = "Timesatser"

# --- Update data on the "Portefølje_F2026" sheet ---
$ws = $wb.Worksheets.Item("Portefølje_F2026")

# Row 2: Søren Erbs Poulsen (SOEB)
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 408
$ws.Range("E2").Value = 50
$ws.Range("G2").Value = 658
$ws.Range("H2").Value = 0

# Row 4: Rune Kier Nielsen (RUNI)
$ws.Range("D4").Value = 178
$ws.Range("E4").Value = 0

# Row 8: Marton Major (MMAJ)
$ws.Range("B8").Value = 110
$ws.Range("C8").Value = 437
$ws.Range("D8").Value = 52
$ws.Range("E8").Value = 59

# Row 9: NN (ufordelt)
$ws.Range("B9").Value = 234
$ws.Range("G9").Value = 327

# Row 10 label: drop the "*** " prefix
$ws.Range("A10").Value = "Projektbudget [kr]"

# Row 11: label and relabeled totals
$ws.Range("A11").Value = "Projektomkostning [kr]"
$ws.Range("B11").Value = 370036
$ws.Range("C11").Value = 231177
$ws.Range("D11").Value = 399181
$ws.Range("E11").Value = 299786
